$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the reference-[5] paragraph (last body paragraph before the sectPr).
# ---------------------------------------------------------------------------
$p5 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*5]*陆雄文*") {
        $p5 = $cand
    }
}

# ---------------------------------------------------------------------------
# Step 1: insert the rebuilt "[5] ... " paragraph (no pPr, the space+name runs
# merged into one run) right after the existing paragraph, via InsertXML at a
# point collapsed just before the paragraph's own mark -- InsertXML there
# inserts a brand-new, cleanly-built paragraph *after* the current one without
# disturbing the current paragraph's own content.
# ---------------------------------------------------------------------------
$insertAfterP5 = $d.Range($p5.Range.End - 1, $p5.Range.End - 1)

$xmlP5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>[</w:t></w:r><w:r><w:t>5]</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 陆雄文</w:t></w:r><w:r><w:t>.管理学大辞典:[M].上海:上海世纪出版股份有限公司上海辞书出版社,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertAfterP5.InsertXML($xmlP5)

# ---------------------------------------------------------------------------
# Step 2: remove the *old* paragraph's run content (keep its end-of-paragraph
# mark for now), then delete that now-empty mark so it merges forward into
# the freshly inserted paragraph -- a paragraph-mark delete keeps the
# *second* paragraph's own (pPr-less) mark, which is exactly how the old
# "[hint=eastAsia]" pPr gets dropped.
# ---------------------------------------------------------------------------
$oldContent = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$oldContent.Delete()

$oldMark = $d.Range($p5.Range.End - 1, $p5.Range.End)
$oldMark.Delete()

# ---------------------------------------------------------------------------
# Re-locate the rebuilt "[5]" paragraph, now pPr-free.
# ---------------------------------------------------------------------------
$p5new = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*5]*陆雄文*") {
        $p5new = $cand
    }
}

# ---------------------------------------------------------------------------
# Step 3: insert the new reference-[6] paragraph right after it, carrying its
# own pPr (spacing + rPr) and the three runs ("[", "6] ", and the formatted
# citation text).
# ---------------------------------------------------------------------------
$insertAfterP5new = $d.Range($p5new.Range.End - 1, $p5new.Range.End - 1)

$xmlP6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="300" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>[</w:t></w:r><w:r><w:t xml:space="preserve">6] </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>陶树平.数据库原理与运用[M].科学出版社，2005.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertAfterP5new.InsertXML($xmlP6)

# ---------------------------------------------------------------------------
# Step 4: append one more, completely empty paragraph ("<w:p/>") at the very
# end of the document body. InsertXML refuses to materialise a paragraph that
# carries no runs at all, so insert one with a throw-away run of text and
# then delete just that text, leaving a bare paragraph mark behind.
# ---------------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xmlEmpty = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>X</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endRange.InsertXML($xmlEmpty)

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$placeholder = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$placeholder.Delete()

Write-Output "done"
